$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values - force text to avoid Excel numeric coercion
$dUpdates = @{
    "D2" = '66.405.62'
    "D3" = '3.301.69'
    "D5" = '570.83'
    "D6" = '181.67'
    "D10" = '6.61'
    "D12" = '3.875.60'
    "D14" = '27.08'
    "D15" = '66.492.38'
    "D17" = '3.291.65'
    "D18" = '13.63'
    "D20" = '429.14'
    "D21" = '7.59'
    "D22" = '73.48'
    "D25" = '0.0000117'
    "D30" = '22.70'
    "D31" = '5.30'
    "D34" = '6.75'
    "D35" = '1.49'
    "D36" = '159.77'
    "D38" = '27.08'
    "D39" = '2.794.57'
    "D40" = '0.787'
    "D41" = '4.42'
    "D42" = '6.16'
    "D43" = '0.0672'
    "D44" = '40.13'
    "D45" = '24.28'
    "D46" = '2.34'
    "D47" = '319.76'
    "D49" = '0.980'
    "D50" = '6.15'
}
foreach ($addr in $dUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $dUpdates[$addr]
    $cell.ClearFormats()
}

# Update Volume(1h) (column E) values
$eUpdates = @{
    "E2" = '  -4.59%  '
    "E3" = '  -1.44%  '
    "E5" = '  -3.62%  '
    "E6" = '  -6.10%  '
    "E7" = '  +0.08%  '
    "E8" = '  -1.50%  '
    "E9" = '  -3.92%  '
    "E10" = '  -2.24%  '
    "E11" = '  -4.90%  '
    "E12" = '  -1.40%  '
    "E13" = '  -0.89%  '
    "E14" = '  -4.73%  '
    "E15" = '  -4.47%  '
    "E16" = '  -2.99%  '
    "E17" = '  -1.55%  '
    "E18" = '  -0.74%  '
    "E19" = '  -2.83%  '
    "E20" = '  -2.76%  '
    "E21" = '  -2.45%  '
    "E22" = '  -0.13%  '
    "E23" = '  -0.09%  '
    "E25" = '  -3.56%  '
    "E26" = '  +0.24%  '
    "E27" = '  -5.99%  '
    "E28" = '  -1.36%  '
    "E29" = '  -2.29%  '
    "E30" = '  -1.69%  '
    "E31" = '  -5.37%  '
    "E32" = '  +0.00%  '
    "E33" = '  -3.66%  '
    "E34" = '  -4.13%  '
    "E35" = '  -1.42%  '
    "E36" = '  -2.81%  '
    "E37" = '  -4.43%  '
    "E38" = '  -0.35%  '
    "E39" = '  +1.04%  '
    "E40" = '  -3.12%  '
    "E41" = '  -3.51%  '
    "E42" = '  -5.18%  '
    "E43" = '  -2.41%  '
    "E44" = '  -1.23%  '
    "E45" = '  -4.62%  '
    "E46" = '  -7.03%  '
    "E47" = '  -6.95%  '
    "E48" = '  -4.27%  '
    "E49" = '  -2.97%  '
    "E50" = '  -2.10%  '
    "E51" = '  -1.28%  '
}
foreach ($addr in $eUpdates.Keys) {
    $ws.Range($addr).Value = $eUpdates[$addr]
}
